# Apply cryptocurrency price/volume updates to sheet1 (cryptos.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'35.372.64"
$ws.Range("E2").Value = "  +0.52%  "

$ws.Range("D3").Value = "'1.890.02"
$ws.Range("E3").Value = "  -0.53%  "

$ws.Range("E4").Value = "  -0.70%  "

$ws.Range("D5").Value = "'246.53"
$ws.Range("E5").Value = "  -2.34%  "

$ws.Range("D6").Value = "'0.691"
$ws.Range("E6").Value = "  -0.49%  "

$ws.Range("E7").Value = "  -0.80%  "

$ws.Range("D8").Value = "'43.27"
$ws.Range("E8").Value = "  +5.19%  "

$ws.Range("E9").Value = "  -1.73%  "

$ws.Range("D10").Value = "'53.87"
$ws.Range("E10").Value = "  +1.71%  "

$ws.Range("D11").Value = "'0.0743"
$ws.Range("E11").Value = "  -1.98%  "

$ws.Range("D12").Value = "'0.0970"
$ws.Range("E12").Value = "  -1.48%  "

$ws.Range("D13").Value = "'13.32"
$ws.Range("E13").Value = "  +1.90%  "

$ws.Range("D14").Value = "'2.165.90"
$ws.Range("E14").Value = "  -0.49%  "

$ws.Range("D15").Value = "'0.759"
$ws.Range("E15").Value = "  +2.90%  "

$ws.Range("E16").Value = "  -1.30%  "

$ws.Range("D17").Value = "'1.872.02"
$ws.Range("E17").Value = "  -1.77%  "

$ws.Range("D18").Value = "'35.394.23"
$ws.Range("E18").Value = "  +0.61%  "

$ws.Range("D19").Value = "'73.22"
$ws.Range("E19").Value = "  -0.79%  "

$ws.Range("D20").Value = "'0.0₃0825"
$ws.Range("E20").Value = "  -1.28%  "

$ws.Range("D21").Value = "'245.04"
$ws.Range("E21").Value = "  +1.02%  "

$ws.Range("D22").Value = "'12.81"
$ws.Range("E22").Value = "  -1.04%  "

$ws.Range("E23").Value = "  -1.84%  "

$ws.Range("E24").Value = "  +11.13%  "

$ws.Range("E25").Value = "  -0.83%  "

$ws.Range("D26").Value = "'2.14"
$ws.Range("E26").Value = "  -7.78%  "

$ws.Range("D27").Value = "'166.20"
$ws.Range("E27").Value = "  -0.53%  "

$ws.Range("D28").Value = "'8.51"
$ws.Range("E28").Value = "  -1.19%  "

$ws.Range("D29").Value = "'18.36"
$ws.Range("E29").Value = "  -0.75%  "

$ws.Range("D30").Value = "'0.127"
$ws.Range("E30").Value = "  -1.55%  "

$ws.Range("D31").Value = "'4.128.44"
$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("D32").Value = "'1.75"
$ws.Range("E32").Value = "  +10.68%  "

$ws.Range("E33").Value = "  -1.51%  "

$ws.Range("E34").Value = "  -3.78%  "

$ws.Range("E35").Value = "  -0.61%  "

$ws.Range("D36").Value = "'1.87"
$ws.Range("E36").Value = "  -12.26%  "

$ws.Range("E37").Value = "  -0.76%  "

$ws.Range("D38").Value = "'0.847"
$ws.Range("E38").Value = "  -0.91%  "

$ws.Range("D40").Value = "'0.0694"
$ws.Range("E40").Value = "  +7.39%  "

$ws.Range("E41").Value = "  +2.60%  "

$ws.Range("D42").Value = "'17.26"
$ws.Range("E42").Value = "  +1.00%  "

$ws.Range("D43").Value = "'96.97"
$ws.Range("E43").Value = "  -3.01%  "

$ws.Range("E44").Value = "  -2.38%  "

$ws.Range("D45").Value = "'1.299.09"

$ws.Range("E46").Value = "  -4.28%  "

$ws.Range("E48").Value = "  -1.03%  "

$ws.Range("E49").Value = "  +3.58%  "

$ws.Range("E50").Value = "  -0.58%  "

$ws.Range("E51").Value = "  -5.36%  "
